$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format
# first, so Excel keeps storing them as text (matching the original
# inline-string cell type) instead of coercing to a number.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = '65.110.89'
$ws.Range("E2").Value = '  +1.88%  '
$ws.Range("D3").Value = '3.163.82'
$ws.Range("E3").Value = '  +3.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '575.35'
$ws.Range("E5").Value = '  +2.78%  '
$ws.Range("D6").Value = '150.62'
$ws.Range("E6").Value = '  +5.69%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.163.76'
$ws.Range("E8").Value = '  +3.30%  '
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("D11").Value = '6.11'
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = '0.500'
$ws.Range("E12").Value = '  +3.84%  '
$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").Value = '  +14.73%  '
$ws.Range("D14").Value = '37.30'
$ws.Range("E14").Value = '  +5.48%  '
$ws.Range("D15").Value = '3.677.54'
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").Value = '65.088.89'
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("D17").Value = '3.149.23'
$ws.Range("E17").Value = '  +2.88%  '
$ws.Range("D18").Value = '7.13'
$ws.Range("E18").Value = '  +4.72%  '
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '508.47'
$ws.Range("E20").Value = '  +4.42%  '
$ws.Range("D21").Value = '14.91'
$ws.Range("E21").Value = '  +3.79%  '
$ws.Range("D22").Value = '0.721'
$ws.Range("E22").Value = '  +4.25%  '
$ws.Range("D23").Value = '15.36'
$ws.Range("E23").Value = '  +4.24%  '
$ws.Range("D24").Value = '7.76'
$ws.Range("E24").Value = '  +3.12%  '
$ws.Range("D25").Value = '84.51'
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '8.98'
$ws.Range("E27").Value = '  +9.67%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.91'
$ws.Range("E28").Value = '  +3.68%  '
$ws.Range("D29").Value = '2.17'
$ws.Range("E29").Value = '  +5.65%  '
$ws.Range("D30").Value = '2.82'
$ws.Range("E30").Value = '  +9.91%  '
$ws.Range("D31").Value = '27.74'
$ws.Range("E31").Value = '  +4.65%  '
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("E33").Value = '  +2.64%  '
$ws.Range("D34").Value = '6.21'
$ws.Range("E34").Value = '  +7.97%  '
$ws.Range("D35").Value = '6.54'
$ws.Range("E35").Value = '  +4.12%  '
$ws.Range("D36").Value = '54.80'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '0.0904'
$ws.Range("E37").Value = '  +10.56%  '
$ws.Range("D38").Value = '467.47'
$ws.Range("E38").Value = '  +5.69%  '
$ws.Range("D39").Value = '0.0422'
$ws.Range("E39").Value = '  +2.49%  '
$ws.Range("D40").Value = '3.01'
$ws.Range("E40").Value = '  +10.36%  '
$ws.Range("D41").Value = '8.69'
$ws.Range("E41").Value = '  +3.93%  '
$ws.Range("D42").Value = '3.061.29'
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").Value = '0.118'
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("D44").Value = '2.46'
$ws.Range("E44").Value = '  +9.65%  '
$ws.Range("D45").Value = '0.284'
$ws.Range("E45").Value = '  +2.81%  '
$ws.Range("D46").Value = '28.65'
$ws.Range("E46").Value = '  +2.36%  '
$ws.Range("D47").Value = '0.0₃0590'
$ws.Range("E47").Value = '  +13.85%  '
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").Value = '0.115'
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  +5.26%  '
$ws.Range("D51").Value = '119.36'
$ws.Range("E51").Value = '  +1.54%  '
